{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// 1) \"SUPERIOR COURT OF \" -> \"SUPERIOR COURT OF test12\"\n// 2) \"Case No.: \" -> \"Case No.: test1\"\n// 3) \"...located at ,  will...\" -> \"...located at test1,  will...\"\n// 4) first empty paragraph right after \"CONCLUSION\" heading -> \"test1\"\n// 5) next empty paragraph (after \"Respectfully submitted,\" + blank line) -> \"test\"\nfor (let i = 0; i < items.length; i++) {\n  const p = items[i];\n  const text = p.text;\n\n  if (text === \"SUPERIOR COURT OF \") {\n    p.insertText(\"SUPERIOR COURT OF test12\", \"Replace\");\n  } else if (text === \"Case No.: \") {\n    p.insertText(\"Case No.: test1\", \"Replace\");\n  } else if (text.indexOf(\"located at ,\") !== -1) {\n    const updated = text.replace(\"located at ,\", \"located at test1,\");\n    p.insertText(updated, \"Replace\");\n  } else if (text === \"\" && items[i - 1] && items[i - 1].text === \"CONCLUSION\") {\n    p.insertText(\"test1\", \"Replace\");\n  } else if (text === \"\" && items[i + 1] && items[i + 1].text === \"Attorney for \") {\n    p.insertText(\"test\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Rewrite a paragraph's text to $newText while preserving the existing run's\n# xml:space=\"preserve\" whitespace handling as much as possible. Plain\n# `Range.Text = ...` (or a Find/Replace) rewrites the whole run and lets the\n# engine re-derive xml:space purely from the new string's own leading /\n# trailing whitespace, which can silently drop xml:space=\"preserve\" when the\n# inserted text happens to not need it at its own boundaries. Splitting the\n# edit into \"delete the trailing part that changed, then re-append the new\n# trailing part via a collapsed end-of-range InsertAfter\" keeps the original\n# run's preserve-space handling intact, matching how Word itself edits text\n# in place.\nfunction Set-ParagraphText($para, [string]$newText) {\n    $r = $para.Range\n    $oldText = $r.Text\n    if ($oldText.Length -gt 0) {\n        $oldText = $oldText.Substring(0, $oldText.Length - 1)  # drop trailing paragraph mark\n    }\n\n    # Find the longest common prefix between old and new text; only the tail\n    # after that prefix actually changed, so only it needs to be touched.\n    $maxPrefix = [Math]::Min($oldText.Length, $newText.Length)\n    $prefixLen = 0\n    while ($prefixLen -lt $maxPrefix -and $oldText[$prefixLen] -eq $newText[$prefixLen]) {\n        $prefixLen++\n    }\n    $newTail = $newText.Substring($prefixLen)\n\n    if ($prefixLen -lt $oldText.Length) {\n        $tailStart = $r.Start + $prefixLen\n        $tailEnd = $r.Start + $oldText.Length\n        $tailRange = $d.Range($tailStart, $tailEnd)\n        $tailRange.Delete()\n    }\n\n    $r2 = $para.Range\n    if ($r2.End -gt $r2.Start) {\n        $r2.End = $r2.End - 1\n    }\n    $r2.Collapse(0)\n    $r2.InsertAfter($newTail)\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t.Length -gt 0) {\n        $t = $t.Substring(0, $t.Length - 1)\n    }\n\n    if ($t -eq \"SUPERIOR COURT OF \") {\n        Set-ParagraphText $p \"SUPERIOR COURT OF test12\"\n    }\n    elseif ($t -eq \"Case No.: \") {\n        Set-ParagraphText $p \"Case No.: test1\"\n    }\n    elseif ($t -like \"*located at ,*\") {\n        $updated = $t.Replace(\"located at ,\", \"located at test1,\")\n        Set-ParagraphText $p $updated\n    }\n    elseif ($t -eq \"\" -and $i -gt 1 -and ($d.Paragraphs.Item($i - 1).Range.Text.TrimEnd(\"`r\")) -eq \"CONCLUSION\") {\n        Set-ParagraphText $p \"test1\"\n    }\n    elseif ($t -eq \"\" -and $i -lt $count -and ($d.Paragraphs.Item($i + 1).Range.Text.TrimEnd(\"`r\")) -eq \"Attorney for \") {\n        Set-ParagraphText $p \"test\"\n    }\n}\n"}
